# Updates the cryptos price list (columns D = Price, E = Volume(1h)) to the
# refreshed snapshot, and swaps the rank-49/50 rows (ordi <-> MultiversX)
# to match the newly published ordering.
#
# Several "Price" values look like plain decimal numbers (e.g. "1.00",
# "314.71"). Assigning such a string straight to Range.Value lets Excel's
# COM layer auto-coerce it to a numeric cell, which would silently change
# the cell's stored type. Set-TextValue forces the cell to Text format
# just long enough to accept the literal string, then clears the
# now-unneeded formatting so the cell is left with no explicit style
# (matching how the sheet looked before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "41.524.94"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.470.60"
$ws.Range("E3").Value = "  -0.62%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue $ws.Range("D5") "314.71"
$ws.Range("E5").Value = "  +0.54%  "
Set-TextValue $ws.Range("D6") "91.96"
$ws.Range("E6").Value = "  -2.92%  "
Set-TextValue $ws.Range("D7") "0.549"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.18%  "
Set-TextValue $ws.Range("D9") "0.515"
$ws.Range("E9").Value = "  +3.00%  "
Set-TextValue $ws.Range("D10") "32.24"
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "2.850.78"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("E14").Value = "  -2.13%  "
Set-TextValue $ws.Range("D15") "15.99"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").Value = "2.467.83"
$ws.Range("E16").Value = "  +1.45%  "
Set-TextValue $ws.Range("D17") "0.770"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "41.516.07"
$ws.Range("E18").Value = "  +0.25%  "
Set-TextValue $ws.Range("D19") "6.47"
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("E20").Value = "  +2.34%  "
Set-TextValue $ws.Range("D21") "71.30"
$ws.Range("E21").Value = "  +3.36%  "
Set-TextValue $ws.Range("D22") "11.08"
$ws.Range("E22").Value = "  -2.01%  "
Set-TextValue $ws.Range("D23") "235.84"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -0.59%  "
Set-TextValue $ws.Range("D27") "24.65"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("E28").Value = "  -0.60%  "
Set-TextValue $ws.Range("D29") "9.68"
$ws.Range("E29").Value = "  -0.86%  "
Set-TextValue $ws.Range("D30") "35.29"
$ws.Range("E30").Value = "  -3.38%  "
Set-TextValue $ws.Range("D31") "155.77"
$ws.Range("E31").Value = "  +2.53%  "
Set-TextValue $ws.Range("D32") "5.44"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  +0.89%  "
Set-TextValue $ws.Range("D35") "17.21"
$ws.Range("E35").Value = "  -4.40%  "
Set-TextValue $ws.Range("D36") "2.87"
$ws.Range("E36").Value = "  -7.23%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  -0.56%  "
Set-TextValue $ws.Range("D39") "1.78"
$ws.Range("E39").Value = "  -5.28%  "
Set-TextValue $ws.Range("D40") "2.23"
$ws.Range("E40").Value = "  -12.87%  "
Set-TextValue $ws.Range("D41") "4.04"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "1.942.29"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("E44").Value = "  -1.40%  "
Set-TextValue $ws.Range("D45") "18.39"
$ws.Range("E45").Value = "  -6.34%  "
Set-TextValue $ws.Range("D46") "2.93"
$ws.Range("E46").Value = "  -3.42%  "
Set-TextValue $ws.Range("D47") "9.04"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").Value = "2.708.03"
$ws.Range("E48").Value = "  -0.84%  "
Set-TextValue $ws.Range("D49") "96.90"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D50") "66.98"
$ws.Range("E50").Value = "  -4.16%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D51") "52.77"
$ws.Range("E51").Value = "  +3.28%  "
